$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows with Price (D) and Volume(1h) (E) updates ---
$ws.Range("D2").Value = "'58.940.88"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").Value = "'2.637.32"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("D5").Value = "'518.34"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").Value = "'146.17"
$ws.Range("E6").Value = "  -2.02%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("D8").Value = "'0.576"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'2.648.20"
$ws.Range("E9").Value = "  -1.24%  "
$ws.Range("D10").Value = "'6.36"
$ws.Range("E10").Value = "  -3.55%  "
$ws.Range("D11").Value = "'0.105"
$ws.Range("E11").Value = "  -1.12%  "
$ws.Range("D12").Value = "'0.336"
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("D14").Value = "'3.099.08"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "'58.913.08"
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("D16").Value = "'20.89"
$ws.Range("E16").Value = "  -2.96%  "
$ws.Range("D18").Value = "'2.641.17"
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").Value = "'349.31"
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("D20").Value = "'4.48"
$ws.Range("E20").Value = "  -3.55%  "
$ws.Range("D21").Value = "'10.27"
$ws.Range("E21").Value = "  -2.98%  "
$ws.Range("D22").Value = "'6.18"
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").Value = "'61.90"
$ws.Range("E24").Value = "  +1.46%  "
$ws.Range("D26").Value = "'0.164"
$ws.Range("E26").Value = "  +1.41%  "
$ws.Range("D27").Value = "'0.996"
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("D28").Value = "'0.0₃0806"
$ws.Range("E28").Value = "  -4.14%  "
$ws.Range("D29").Value = "'7.07"
$ws.Range("E29").Value = "  -1.66%  "
$ws.Range("D31").Value = "'6.28"
$ws.Range("E31").Value = "  -7.94%  "
$ws.Range("D34").Value = "'149.38"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "'0.972"
$ws.Range("E35").Value = "  -8.06%  "
$ws.Range("D36").Value = "'4.03"
$ws.Range("E36").Value = "  -1.27%  "
$ws.Range("D38").Value = "'36.62"
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("D39").Value = "'0.847"
$ws.Range("E39").Value = "  -3.72%  "
$ws.Range("D40").Value = "'1.42"
$ws.Range("E40").Value = "  -1.50%  "
$ws.Range("D41").Value = "'3.64"
$ws.Range("E41").Value = "  -2.40%  "
$ws.Range("D42").Value = "'279.39"
$ws.Range("E42").Value = "  -1.91%  "
$ws.Range("D44").Value = "'0.0989"
$ws.Range("E44").Value = "  -1.23%  "
$ws.Range("D45").Value = "'19.76"
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("D46").Value = "'0.602"
$ws.Range("E46").Value = "  -4.70%  "
$ws.Range("D50").Value = "'0.0230"
$ws.Range("E50").Value = "  -1.66%  "
$ws.Range("D51").Value = "'4.70"
$ws.Range("E51").Value = "  -2.43%  "

# --- Rows with only Volume(1h) (E) updates ---
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("E17").Value = "  -2.09%  "
$ws.Range("E25").Value = "  -3.02%  "
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("E37").Value = "  -0.83%  "
$ws.Range("E43").Value = "  +0.52%  "

# --- Rows with full Coin/Link/Price/Volume updates (reordering) ---
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'1.58"
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'18.90"
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "'2.055.23"
$ws.Range("E47").Value = "  +3.23%  "
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").Value = "'0.0525"
$ws.Range("E48").Value = "  -4.04%  "
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").Value = "'10.31"
$ws.Range("E49").Value = "  +0.51%  "
